$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.35283127929076841
$ws.Range("A2").Value = -0.0099999997313524602
$ws.Range("A3").Value = -0.0089999997274645693
$ws.Range("A4").Value = 0.061997179486139586
$ws.Range("A5").Value = -0.0059999997313893161
$ws.Range("A6").Value = -0.0059999997204833733
$ws.Range("A7").Value = -0.019999999665028412
$ws.Range("A8").Value = -0.042087514046022179
$ws.Range("A9").Value = -0.005999999715085913
$ws.Range("A10").Value = -0.0059999997141702011
$ws.Range("A11").Value = -0.0044999997199468567
$ws.Range("A12").Value = 0.043680799111356805
$ws.Range("A13").Value = -0.0059999997091404467
$ws.Range("A14").Value = -0.011999999684317864
$ws.Range("A15").Value = -0.0059999997064714705
$ws.Range("A16").Value = -0.0059999997053876708
$ws.Range("A17").Value = -0.0059999997039605901
$ws.Range("A18").Value = -0.0089999996919578606
$ws.Range("A19").Value = -0.0089999997342053994
$ws.Range("A20").Value = -0.0089999997245620023
$ws.Range("A21").Value = -0.076383652242593136
$ws.Range("A22").Value = -0.0089999997219392114
$ws.Range("A23").Value = -0.0089999997248551011
$ws.Range("A24").Value = -0.041999999592031934
$ws.Range("A25").Value = -0.041999999589684478
$ws.Range("A26").Value = -0.0059999997192115018
$ws.Range("A27").Value = -0.0059999997174249309
$ws.Range("A28").Value = -0.0059999997107560432
$ws.Range("A29").Value = -0.011999999682771545
$ws.Range("A30").Value = -0.019999999649385813
$ws.Range("A31").Value = -0.014999999665514352
$ws.Range("A32").Value = -0.020999999641650113
$ws.Range("A33").Value = -0.0059999996997248672
